# Update column G ("K" - strikeouts) values for rows 2-10 and 12
# to reflect regenerated save_data (using K instead of Strike#).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    2  = 3
    3  = 1
    4  = 6
    5  = 3
    6  = 5
    7  = 3
    8  = 1
    9  = 5
    10 = 4
    12 = 1
}

foreach ($row in $values.Keys) {
    $ws.Cells.Item($row, 7).Value = $values[$row]
}
